$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38
$c = $ws.Range('A38'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B38').Value = 'Major League Soccer'
$ws.Range('C38').Value = 'Austin'
$ws.Range('D38').Value = 'San Jose Earthquakes'
$ws.Range('E38').Value = 'Home Win'
$c = $ws.Range('F38'); $c.NumberFormat = '@'; $c.Value = '53.83%'; $c.Style = 'Normal'
$ws.Range('G38').Value = 2.1
$c = $ws.Range('H38'); $c.NumberFormat = '@'; $c.Value = '11.91%'; $c.Style = 'Normal'
$ws.Range('I38').Value = 1.3
$ws.Range('J38').Value = 0.01185935992232205
$ws.Range('K38').Value = 0.1185935992232205
$ws.Range('L38').Value = 'Pending'

# Row 39
$c = $ws.Range('A39'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B39').Value = 'Major League Soccer'
$ws.Range('C39').Value = 'Sporting Kansas City'
$ws.Range('D39').Value = 'Colorado Rapids'
$ws.Range('E39').Value = 'Home Win'
$c = $ws.Range('F39'); $c.NumberFormat = '@'; $c.Value = '58.15%'; $c.Style = 'Normal'
$ws.Range('G39').Value = 2
$c = $ws.Range('H39'); $c.NumberFormat = '@'; $c.Value = '15.14%'; $c.Style = 'Normal'
$ws.Range('I39').Value = 1.7
$ws.Range('J39').Value = 0.01629821757331635
$ws.Range('K39').Value = 0.1629821757331635
$ws.Range('L39').Value = 'Pending'

# Row 40
$c = $ws.Range('A40'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B40').Value = 'Major League Soccer'
$ws.Range('C40').Value = 'Minnesota United FC'
$ws.Range('D40').Value = 'Portland Timbers'
$ws.Range('E40').Value = 'Home Win'
$c = $ws.Range('F40'); $c.NumberFormat = '@'; $c.Value = '73.61%'; $c.Style = 'Normal'
$ws.Range('G40').Value = 1.75
$c = $ws.Range('H40'); $c.NumberFormat = '@'; $c.Value = '27.53%'; $c.Style = 'Normal'
$ws.Range('I40').Value = 4.1
$ws.Range('J40').Value = 0.03843006390250685
$ws.Range('K40').Value = 0.3843006390250685
$ws.Range('L40').Value = 'Pending'

# Row 41
$c = $ws.Range('A41'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B41').Value = 'Major League Soccer'
$ws.Range('C41').Value = 'Nashville SC'
$ws.Range('D41').Value = 'Atlanta United FC'
$ws.Range('E41').Value = 'Home Win'
$c = $ws.Range('F41'); $c.NumberFormat = '@'; $c.Value = '87.87%'; $c.Style = 'Normal'
$ws.Range('G41').Value = 1.44
$c = $ws.Range('H41'); $c.NumberFormat = '@'; $c.Value = '25.27%'; $c.Style = 'Normal'
$ws.Range('I41').Value = 5.3
$ws.Range('J41').Value = 0.05
$ws.Range('K41').Value = 0.6030328745841734
$ws.Range('L41').Value = 'Pending'

# Row 42
$c = $ws.Range('A42'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B42').Value = 'Eredivisie'
$ws.Range('C42').Value = 'Fortuna Sittard'
$ws.Range('D42').Value = 'NEC Nijmegen'
$ws.Range('E42').Value = 'Away Win'
$c = $ws.Range('F42'); $c.NumberFormat = '@'; $c.Value = '49.38%'; $c.Style = 'Normal'
$ws.Range('G42').Value = 2.2
$c = $ws.Range('H42'); $c.NumberFormat = '@'; $c.Value = '7.54%'; $c.Style = 'Normal'
$ws.Range('I42').Value = 0.8
$ws.Range('J42').Value = 0.007187961864934366
$ws.Range('K42').Value = 0.07187961864934365
$ws.Range('L42').Value = 'Pending'

# Row 43
$c = $ws.Range('A43'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B43').Value = 'Jupiler Pro League'
$ws.Range('C43').Value = 'Gent'
$ws.Range('D43').Value = 'Club Brugge KV'
$ws.Range('E43').Value = 'Away Win'
$c = $ws.Range('F43'); $c.NumberFormat = '@'; $c.Value = '80.89%'; $c.Style = 'Normal'
$ws.Range('G43').Value = 1.62
$c = $ws.Range('H43'); $c.NumberFormat = '@'; $c.Value = '29.74%'; $c.Style = 'Normal'
$ws.Range('I43').Value = 5.3
$ws.Range('J43').Value = 0.05
$ws.Range('K43').Value = 0.500760157021969
$ws.Range('L43').Value = 'Pending'

# Row 44
$c = $ws.Range('A44'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B44').Value = 'Eredivisie'
$ws.Range('C44').Value = 'Sparta Rotterdam'
$ws.Range('D44').Value = 'Feyenoord'
$ws.Range('E44').Value = 'Away Win'
$c = $ws.Range('F44'); $c.NumberFormat = '@'; $c.Value = '67.61%'; $c.Style = 'Normal'
$ws.Range('G44').Value = 1.83
$c = $ws.Range('H44'); $c.NumberFormat = '@'; $c.Value = '22.49%'; $c.Style = 'Normal'
$ws.Range('I44').Value = 3
$ws.Range('J44').Value = 0.02858543913583976
$ws.Range('K44').Value = 0.2858543913583976
$ws.Range('L44').Value = 'Pending'

# Row 45
$c = $ws.Range('A45'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B45').Value = 'Eredivisie'
$ws.Range('C45').Value = 'PEC Zwolle'
$ws.Range('D45').Value = 'Utrecht'
$ws.Range('E45').Value = 'Away Win'
$c = $ws.Range('F45'); $c.NumberFormat = '@'; $c.Value = '50.65%'; $c.Style = 'Normal'
$ws.Range('G45').Value = 2.2
$c = $ws.Range('H45'); $c.NumberFormat = '@'; $c.Value = '10.31%'; $c.Style = 'Normal'
$ws.Range('I45').Value = 1
$ws.Range('J45').Value = 0.009524246937548303
$ws.Range('K45').Value = 0.09524246937548302
$ws.Range('L45').Value = 'Pending'

# Row 46
$c = $ws.Range('A46'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B46').Value = 'Premier League'
$ws.Range('C46').Value = 'Brighton'
$ws.Range('D46').Value = 'Manchester City'
$ws.Range('E46').Value = 'Away Win'
$c = $ws.Range('F46'); $c.NumberFormat = '@'; $c.Value = '67.98%'; $c.Style = 'Normal'
$ws.Range('G46').Value = 1.85
$c = $ws.Range('H46'); $c.NumberFormat = '@'; $c.Value = '24.51%'; $c.Style = 'Normal'
$ws.Range('I46').Value = 3.2
$ws.Range('J46').Value = 0.03031818664193643
$ws.Range('K46').Value = 0.3031818664193643
$ws.Range('L46').Value = 'Pending'

# Row 47
$c = $ws.Range('A47'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B47').Value = 'Premier League'
$ws.Range('C47').Value = 'Nottingham Forest'
$ws.Range('D47').Value = 'West Ham'
$ws.Range('E47').Value = 'Home Win'
$c = $ws.Range('F47'); $c.NumberFormat = '@'; $c.Value = '75.77%'; $c.Style = 'Normal'
$ws.Range('G47').Value = 1.67
$c = $ws.Range('H47'); $c.NumberFormat = '@'; $c.Value = '25.27%'; $c.Style = 'Normal'
$ws.Range('I47').Value = 4.2
$ws.Range('J47').Value = 0.03959937087708185
$ws.Range('K47').Value = 0.3959937087708185
$ws.Range('L47').Value = 'Pending'

# Row 48
$c = $ws.Range('A48'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B48').Value = 'Ligue 1'
$ws.Range('C48').Value = 'Angers'
$ws.Range('D48').Value = 'Rennes'
$ws.Range('E48').Value = 'Away Win'
$c = $ws.Range('F48'); $c.NumberFormat = '@'; $c.Value = '61.15%'; $c.Style = 'Normal'
$ws.Range('G48').Value = 2
$c = $ws.Range('H48'); $c.NumberFormat = '@'; $c.Value = '21.08%'; $c.Style = 'Normal'
$ws.Range('I48').Value = 2.4
$ws.Range('J48').Value = 0.02230459415956969
$ws.Range('K48').Value = 0.2230459415956969
$ws.Range('L48').Value = 'Pending'

# Row 49
$c = $ws.Range('A49'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B49').Value = 'Bundesliga'
$ws.Range('C49').Value = 'VfL Wolfsburg'
$ws.Range('D49').Value = 'FSV Mainz 05'
$ws.Range('E49').Value = 'Home Win'
$c = $ws.Range('F49'); $c.NumberFormat = '@'; $c.Value = '61.05%'; $c.Style = 'Normal'
$ws.Range('G49').Value = 2
$c = $ws.Range('H49'); $c.NumberFormat = '@'; $c.Value = '20.89%'; $c.Style = 'Normal'
$ws.Range('I49').Value = 2.3
$ws.Range('J49').Value = 0.0221066358449753
$ws.Range('K49').Value = 0.221066358449753
$ws.Range('L49').Value = 'Pending'

# Row 50
$c = $ws.Range('A50'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B50').Value = 'Jupiler Pro League'
$ws.Range('C50').Value = 'OH Leuven'
$ws.Range('D50').Value = 'Standard Liege'
$ws.Range('E50').Value = 'Home Win'
$c = $ws.Range('F50'); $c.NumberFormat = '@'; $c.Value = '54.09%'; $c.Style = 'Normal'
$ws.Range('G50').Value = 2.15
$c = $ws.Range('H50'); $c.NumberFormat = '@'; $c.Value = '15.14%'; $c.Style = 'Normal'
$ws.Range('I50').Value = 1.5
$ws.Range('J50').Value = 0.0141775166322042
$ws.Range('K50').Value = 0.141775166322042
$ws.Range('L50').Value = 'Pending'

# Row 51
$c = $ws.Range('A51'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B51').Value = 'Eredivisie'
$ws.Range('C51').Value = 'NAC Breda'
$ws.Range('D51').Value = 'AZ Alkmaar'
$ws.Range('E51').Value = 'Away Win'
$c = $ws.Range('F51'); $c.NumberFormat = '@'; $c.Value = '82.19%'; $c.Style = 'Normal'
$ws.Range('G51').Value = 1.57
$c = $ws.Range('H51'); $c.NumberFormat = '@'; $c.Value = '27.74%'; $c.Style = 'Normal'
$ws.Range('I51').Value = 5.3
$ws.Range('J51').Value = 0.05
$ws.Range('K51').Value = 0.5093392164581518
$ws.Range('L51').Value = 'Pending'

# Row 52
$c = $ws.Range('A52'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B52').Value = 'Ligue 1'
$ws.Range('C52').Value = 'Le Havre'
$ws.Range('D52').Value = 'Nice'
$ws.Range('E52').Value = 'Away Win'
$c = $ws.Range('F52'); $c.NumberFormat = '@'; $c.Value = '59.07%'; $c.Style = 'Normal'
$ws.Range('G52').Value = 2.05
$c = $ws.Range('H52'); $c.NumberFormat = '@'; $c.Value = '19.89%'; $c.Style = 'Normal'
$ws.Range('I52').Value = 2.1
$ws.Range('J52').Value = 0.02009732743400738
$ws.Range('K52').Value = 0.2009732743400738
$ws.Range('L52').Value = 'Pending'

# Row 53
$c = $ws.Range('A53'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B53').Value = 'Ligue 1'
$ws.Range('C53').Value = 'Monaco'
$ws.Range('D53').Value = 'Strasbourg'
$ws.Range('E53').Value = 'Home Win'
$c = $ws.Range('F53'); $c.NumberFormat = '@'; $c.Value = '82.37%'; $c.Style = 'Normal'
$ws.Range('G53').Value = 1.55
$c = $ws.Range('H53'); $c.NumberFormat = '@'; $c.Value = '26.40%'; $c.Style = 'Normal'
$ws.Range('I53').Value = 5.3
$ws.Range('J53').Value = 0.05
$ws.Range('K53').Value = 0.5031254641766906
$ws.Range('L53').Value = 'Pending'

# Row 54
$c = $ws.Range('A54'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B54').Value = 'Ligue 1'
$ws.Range('C54').Value = 'Paris FC'
$ws.Range('D54').Value = 'Metz'
$ws.Range('E54').Value = 'Home Win'
$c = $ws.Range('F54'); $c.NumberFormat = '@'; $c.Value = '72.47%'; $c.Style = 'Normal'
$ws.Range('G54').Value = 1.75
$c = $ws.Range('H54'); $c.NumberFormat = '@'; $c.Value = '25.55%'; $c.Style = 'Normal'
$ws.Range('I54').Value = 3.8
$ws.Range('J54').Value = 0.03575920729919848
$ws.Range('K54').Value = 0.3575920729919848
$ws.Range('L54').Value = 'Pending'

# Row 55
$c = $ws.Range('A55'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B55').Value = 'Bundesliga'
$ws.Range('C55').Value = 'Borussia Dortmund'
$ws.Range('D55').Value = 'Union Berlin'
$ws.Range('E55').Value = 'Home Win'
$c = $ws.Range('F55'); $c.NumberFormat = '@'; $c.Value = '90.40%'; $c.Style = 'Normal'
$ws.Range('G55').Value = 1.4
$c = $ws.Range('H55'); $c.NumberFormat = '@'; $c.Value = '25.29%'; $c.Style = 'Normal'
$ws.Range('I55').Value = 5.3
$ws.Range('J55').Value = 0.05
$ws.Range('K55').Value = 0.6639485662539361
$ws.Range('L55').Value = 'Pending'

# Row 56
$c = $ws.Range('A56'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B56').Value = 'Süper Lig'
$ws.Range('C56').Value = 'Istanbul Basaksehir'
$ws.Range('D56').Value = 'Eyüpspor'
$ws.Range('E56').Value = 'Home Win'
$c = $ws.Range('F56'); $c.NumberFormat = '@'; $c.Value = '63.78%'; $c.Style = 'Normal'
$ws.Range('G56').Value = 1.95
$c = $ws.Range('H56'); $c.NumberFormat = '@'; $c.Value = '23.14%'; $c.Style = 'Normal'
$ws.Range('I56').Value = 2.7
$ws.Range('J56').Value = 0.02566256111386165
$ws.Range('K56').Value = 0.2566256111386165
$ws.Range('L56').Value = 'Pending'

# Row 57
$c = $ws.Range('A57'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B57').Value = 'Serie A'
$ws.Range('C57').Value = 'Genoa'
$ws.Range('D57').Value = 'Juventus'
$ws.Range('E57').Value = 'Away Win'
$c = $ws.Range('F57'); $c.NumberFormat = '@'; $c.Value = '71.25%'; $c.Style = 'Normal'
$ws.Range('G57').Value = 1.83
$c = $ws.Range('H57'); $c.NumberFormat = '@'; $c.Value = '29.09%'; $c.Style = 'Normal'
$ws.Range('I57').Value = 3.9
$ws.Range('J57').Value = 0.0366167279866688
$ws.Range('K57').Value = 0.366167279866688
$ws.Range('L57').Value = 'Pending'

# Row 58
$c = $ws.Range('A58'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B58').Value = 'Serie A'
$ws.Range('C58').Value = 'Torino'
$ws.Range('D58').Value = 'Fiorentina'
$ws.Range('E58').Value = 'Away Win'
$c = $ws.Range('F58'); $c.NumberFormat = '@'; $c.Value = '47.32%'; $c.Style = 'Normal'
$ws.Range('G58').Value = 2.4
$c = $ws.Range('H58'); $c.NumberFormat = '@'; $c.Value = '12.43%'; $c.Style = 'Normal'
$ws.Range('I58').Value = 1
$ws.Range('J58').Value = 0.00969258356476905
$ws.Range('K58').Value = 0.0969258356476905
$ws.Range('L58').Value = 'Pending'

# Row 59
$c = $ws.Range('A59'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B59').Value = 'Jupiler Pro League'
$ws.Range('C59').Value = 'Union St. Gilloise'
$ws.Range('D59').Value = 'Anderlecht'
$ws.Range('E59').Value = 'Home Win'
$c = $ws.Range('F59'); $c.NumberFormat = '@'; $c.Value = '80.89%'; $c.Style = 'Normal'
$ws.Range('G59').Value = 1.55
$c = $ws.Range('H59'); $c.NumberFormat = '@'; $c.Value = '24.13%'; $c.Style = 'Normal'
$ws.Range('I59').Value = 4.9
$ws.Range('J59').Value = 0.04615507774293427
$ws.Range('K59').Value = 0.4615507774293426
$ws.Range('L59').Value = 'Pending'

# Row 60
$c = $ws.Range('A60'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B60').Value = 'La Liga'
$ws.Range('C60').Value = 'Real Betis'
$ws.Range('D60').Value = 'Athletic Club'
$ws.Range('E60').Value = 'Away Win'
$c = $ws.Range('F60'); $c.NumberFormat = '@'; $c.Value = '45.39%'; $c.Style = 'Normal'
$ws.Range('G60').Value = 2.45
$c = $ws.Range('H60'); $c.NumberFormat = '@'; $c.Value = '10.08%'; $c.Style = 'Normal'
$ws.Range('I60').Value = 0.8
$ws.Range('J60').Value = 0.007721574202561839
$ws.Range('K60').Value = 0.07721574202561839
$ws.Range('L60').Value = 'Pending'

# Row 61
$c = $ws.Range('A61'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B61').Value = 'Premier League'
$ws.Range('C61').Value = 'Aston Villa'
$ws.Range('D61').Value = 'Crystal Palace'
$ws.Range('E61').Value = 'Home Win'
$c = $ws.Range('F61'); $c.NumberFormat = '@'; $c.Value = '69.61%'; $c.Style = 'Normal'
$ws.Range('G61').Value = 1.83
$c = $ws.Range('H61'); $c.NumberFormat = '@'; $c.Value = '26.11%'; $c.Style = 'Normal'
$ws.Range('I61').Value = 3.5
$ws.Range('J61').Value = 0.03299757013756545
$ws.Range('K61').Value = 0.3299757013756545
$ws.Range('L61').Value = 'Pending'

# Row 62
$c = $ws.Range('A62'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B62').Value = 'Liga de Expansión MX'
$ws.Range('C62').Value = 'Leones Negros UDG'
$ws.Range('D62').Value = 'Tlaxcala'
$ws.Range('E62').Value = 'Home Win'
$c = $ws.Range('F62'); $c.NumberFormat = '@'; $c.Value = '84.10%'; $c.Style = 'Normal'
$ws.Range('G62').Value = 1.5
$c = $ws.Range('H62'); $c.NumberFormat = '@'; $c.Value = '24.89%'; $c.Style = 'Normal'
$ws.Range('I62').Value = 5.3
$ws.Range('J62').Value = 0.05
$ws.Range('K62').Value = 0.5230112549280651
$ws.Range('L62').Value = 'Pending'

# Row 63
$c = $ws.Range('A63'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B63').Value = 'Süper Lig'
$ws.Range('C63').Value = 'Alanyaspor'
$ws.Range('D63').Value = 'Besiktas'
$ws.Range('E63').Value = 'Away Win'
$c = $ws.Range('F63'); $c.NumberFormat = '@'; $c.Value = '54.30%'; $c.Style = 'Normal'
$ws.Range('G63').Value = 2.1
$c = $ws.Range('H63'); $c.NumberFormat = '@'; $c.Value = '12.90%'; $c.Style = 'Normal'
$ws.Range('I63').Value = 1.4
$ws.Range('J63').Value = 0.01276081193961335
$ws.Range('K63').Value = 0.1276081193961335
$ws.Range('L63').Value = 'Pending'

# Row 64
$c = $ws.Range('A64'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B64').Value = 'Süper Lig'
$ws.Range('C64').Value = 'Trabzonspor'
$ws.Range('D64').Value = 'Samsunspor'
$ws.Range('E64').Value = 'Home Win'
$c = $ws.Range('F64'); $c.NumberFormat = '@'; $c.Value = '73.63%'; $c.Style = 'Normal'
$ws.Range('G64').Value = 1.73
$c = $ws.Range('H64'); $c.NumberFormat = '@'; $c.Value = '26.11%'; $c.Style = 'Normal'
$ws.Range('I64').Value = 4
$ws.Range('J64').Value = 0.0375062794028851
$ws.Range('K64').Value = 0.375062794028851
$ws.Range('L64').Value = 'Pending'

# Row 65
$c = $ws.Range('A65'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B65').Value = 'Serie A'
$ws.Range('C65').Value = 'Lazio'
$ws.Range('D65').Value = 'Verona'
$ws.Range('E65').Value = 'Home Win'
$c = $ws.Range('F65'); $c.NumberFormat = '@'; $c.Value = '84.63%'; $c.Style = 'Normal'
$ws.Range('G65').Value = 1.5
$c = $ws.Range('H65'); $c.NumberFormat = '@'; $c.Value = '25.67%'; $c.Style = 'Normal'
$ws.Range('I65').Value = 5.3
$ws.Range('J65').Value = 0.05
$ws.Range('K65').Value = 0.5387770815986485
$ws.Range('L65').Value = 'Pending'

# Row 66
$c = $ws.Range('A66'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B66').Value = 'Ligue 1'
$ws.Range('C66').Value = 'Lyon'
$ws.Range('D66').Value = 'Marseille'
$ws.Range('E66').Value = 'Away Win'
$c = $ws.Range('F66'); $c.NumberFormat = '@'; $c.Value = '45.72%'; $c.Style = 'Normal'
$ws.Range('G66').Value = 2.3
$c = $ws.Range('H66'); $c.NumberFormat = '@'; $c.Value = '4.10%'; $c.Style = 'Normal'
$ws.Range('I66').Value = 0.4
$ws.Range('J66').Value = 0.003965179902846045
$ws.Range('K66').Value = 0.03965179902846044
$ws.Range('L66').Value = 'Pending'

# Row 67
$c = $ws.Range('A67'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B67').Value = 'La Liga'
$ws.Range('C67').Value = 'Rayo Vallecano'
$ws.Range('D67').Value = 'Barcelona'
$ws.Range('E67').Value = 'Away Win'
$c = $ws.Range('F67'); $c.NumberFormat = '@'; $c.Value = '91.11%'; $c.Style = 'Normal'
$ws.Range('G67').Value = 1.4
$c = $ws.Range('H67'); $c.NumberFormat = '@'; $c.Value = '26.29%'; $c.Style = 'Normal'
$ws.Range('I67').Value = 5.3
$ws.Range('J67').Value = 0.05
$ws.Range('K67').Value = 0.6890210134428515
$ws.Range('L67').Value = 'Pending'

# Row 68
$c = $ws.Range('A68'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B68').Value = 'Primeira Liga'
$ws.Range('C68').Value = 'Rio Ave'
$ws.Range('D68').Value = 'SC Braga'
$ws.Range('E68').Value = 'Away Win'
$c = $ws.Range('F68'); $c.NumberFormat = '@'; $c.Value = '78.19%'; $c.Style = 'Normal'
$ws.Range('G68').Value = 1.67
$c = $ws.Range('H68'); $c.NumberFormat = '@'; $c.Value = '29.28%'; $c.Style = 'Normal'
$ws.Range('I68').Value = 4.8
$ws.Range('J68').Value = 0.04564856690064723
$ws.Range('K68').Value = 0.4564856690064723
$ws.Range('L68').Value = 'Pending'

# Row 69
$c = $ws.Range('A69'); $c.NumberFormat = '@'; $c.Value = '2025-08-30'; $c.Style = 'Normal'
$ws.Range('B69').Value = 'Primeira Liga'
$ws.Range('C69').Value = 'Santa Clara'
$ws.Range('D69').Value = 'Estrela'
$ws.Range('E69').Value = 'Home Win'
$c = $ws.Range('F69'); $c.NumberFormat = '@'; $c.Value = '66.74%'; $c.Style = 'Normal'
$ws.Range('G69').Value = 1.83
$c = $ws.Range('H69'); $c.NumberFormat = '@'; $c.Value = '20.91%'; $c.Style = 'Normal'
$ws.Range('I69').Value = 2.8
$ws.Range('J69').Value = 0.02666366764812174
$ws.Range('K69').Value = 0.2666366764812174
$ws.Range('L69').Value = 'Pending'
